# Add files via upload
#
# requirements.docx was edited:
#   - "requirements.txt" -> "streamlit" (1st paragraph)
#   - the spell-check proofErr wrappers around "numpy", "reportlab" and
#     "joblib" are removed (the words are no longer flagged as possible
#     misspellings once the doc has been re-saved from Word)

$d = $word.ActiveDocument

# 1) Simple text swap for the first line.
$d.Content.Find.Execute("requirements.txt", $true, $false, $false, $false, $false,
                         $true, 1, $false, "streamlit", 2)

# 2) Drop the <w:proofErr .../> spell-check markers around a handful of
#    words. There is no direct OM call to toggle proofErr, but deleting a
#    paragraph (including its end-of-paragraph mark) and typing its text
#    back in as a brand-new paragraph reproduces what Word does when it
#    re-flows text that no longer carries a stale proofing flag: the new
#    paragraph/run is created fresh, with no proofErr markers.
function Remove-ProofErrParagraph($wordApp, [string]$text) {
    $doc = $wordApp.ActiveDocument

    $targetIdx = -1
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs($i)
        if ($para.Range.Text.TrimEnd("`r", "`a") -eq $text) {
            $targetIdx = $i
            break
        }
    }
    if ($targetIdx -eq -1) { return }

    $target = $doc.Paragraphs($targetIdx)
    $rng = $target.Range
    $full = $doc.Range($rng.Start, $rng.End)
    $full.Delete()

    if ($targetIdx -gt 1) {
        $prev = $doc.Paragraphs($targetIdx - 1)
        $prev.Range.InsertParagraphAfter()
        $newPara = $doc.Paragraphs($targetIdx)
        $newPara.Range.InsertBefore($text)
    } else {
        $doc.Range(0, 0).InsertParagraphAfter()
        $doc.Paragraphs(1).Range.InsertBefore($text)
    }
}

Remove-ProofErrParagraph $word "numpy"
Remove-ProofErrParagraph $word "reportlab"
Remove-ProofErrParagraph $word "joblib"
